$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 611, shifting rows 611:652 down to 612:653.
$ws.Rows("611:611").Insert()

# Populate the newly inserted row 611 with the new data.
# Force the date-like text into column A as plain text (Excel would
# otherwise auto-convert "2026/01/11" into a date value), then drop the
# cell back to the default "Normal" style so no stray style index lingers
# on the cell (matching the rest of the date column).
$ws.Range("A611").NumberFormat = "@"
$ws.Range("A611").Value = "2026/01/11"
$ws.Range("A611").Style = "Normal"
$ws.Range("B611").Value = "日"
$ws.Range("C611").Value = 5
$ws.Range("D611").Value = 201
